# Finished the "Zth" vs power comparison task, so remove the whole
# to-do bullet paragraph ("Add comparison between Zth and power at
# various times") from under the "Code additions" heading.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Add comparison between Zth and power at various times")) {
        $p.Range.Delete()
    }
}
